$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.704960018034767
$ws.Range("C2").Value = 0.780310385878786
$ws.Range("D2").Value = 0.767261516731427
$ws.Range("E2").Value = 0.757541746745335
$ws.Range("F2").Value = 0.874984008964506
$ws.Range("G2").Value = 0.847206298409435
$ws.Range("H2").Value = 0.842339547896952
$ws.Range("I2").Value = 0.65799487387268
$ws.Range("J2").Value = 0.656332785129309
$ws.Range("K2").Value = 0.703874631903231
$ws.Range("L2").Value = 0.777883926828007
$ws.Range("M2").Value = 0.92816201896394
$ws.Range("N2").Value = 0.562406199574745
